$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data: years 2010..2022 with corresponding B (foreign exchange reserves)
# and C (gold reserves) values, occupying rows 2..14.
$years = @(2010,2011,2012,2013,2014,2015,2016,2017,2018,2019,2020,2021,2022)
$bvals = @(28473.38,31811.48,33115.89,38213.15,38430.18,33303.62,30105.17,31399.49,30727.12,31079.24,32165.22,32501.66,31276.91)
$cvals = @(3389,3389,3389,3389,3389,5666,5924,5924,5956,6264,6264,6264,6464)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = "$($years[$i])年"
    $ws.Cells.Item($row, 2).Value = $bvals[$i]
    $ws.Cells.Item($row, 3).Value = $cvals[$i]
}

# Remove now-unused rows 15..22 (old data extended further down before).
$ws.Range("A15:C22").Clear()
